$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 84, shifting the existing
# rows 84-88 down to 85-89 (matches the dimension growing to A1:R89).
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new weekly price record.
$ws.Range("A84").Value = 6
$ws.Range("B84").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C84").Value = "Metropolitana"
$ws.Range("D84").Value = 44753
$ws.Range("E84").Value = 13
$ws.Range("F84").Value = 100114007
$ws.Range("G84").Value = "Jengibre"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 400
$ws.Range("K84").Value = 13000
$ws.Range("L84").Value = 14000
$ws.Range("M84").Value = 13425
$ws.Range("N84").Value = "$/caja 15 kilos"
$ws.Range("O84").Value = "Perú"
$ws.Range("P84").Value = 895
$ws.Range("Q84").Value = 15
$ws.Range("R84").Value = "Hortaliza"
